$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("out_vars")

$row = 51

# Column A holds text-formatted dates (e.g. "2020-07-20") stored as shared
# strings in the original sheet. A plain .Value assignment of a date-like
# string gets auto-parsed into an Excel date serial, so force the cell to
# Text format before writing it, then restore the cell's style so it matches
# the unstyled data cells around it (no explicit style index).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2020-07-20"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 349396
$ws.Cells.Item($row, 3).Value = 399443
$ws.Cells.Item($row, 4).Value = 79112
$ws.Cells.Item($row, 5).Value = 39485
$ws.Cells.Item($row, 6).Value = 28.43
